# Refresh the crypto price / 1h-volume-change table with the latest
# scrape (GitHub Actions cron job).
#
# Notes:
#  - Price values in column D are stored as TEXT (e.g. "29.195.32",
#    "0.9991"), not numbers -- some contain two decimal points, others
#    are percent-free decimals that Excel would otherwise silently coerce
#    to a real number on assignment. To keep every Price cell as text we
#    set NumberFormat to "@" (Text) just before writing any D-column
#    value that Excel could parse as a number.
#  - Two coins (EnergySwap / BabyDogeCoin) swapped rank this refresh, so
#    rows 48 and 49 get new Coin/Link/Price/Volume values rather than a
#    simple in-place number update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.195.32'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.843.52'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.67'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6877'
$ws.Range('E6').Value = '  -2.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07477'
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.27'
$ws.Range('E10').Value = '  -1.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07662'
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = '1.846.29'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.069'
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '87.47'
$ws.Range('E15').Value = '  -6.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.186'
$ws.Range('D17').Value = '29.194.68'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008181'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').Value = '2.083.79'
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.65'
$ws.Range('E20').Value = '  -5.20%  '
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9994'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.408'
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9996'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1452'
$ws.Range('E25').Value = '  -4.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.59'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.786'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.515'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.286'
$ws.Range('E30').Value = '  +1.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.149'
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.200'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('E33').Value = '  +2.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7610'
$ws.Range('E34').Value = '  -3.74%  '
$ws.Range('E35').Value = '  -2.04%  '
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.682'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').Value = '1.306.43'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01838'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9359'
$ws.Range('E41').Value = '  -2.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.973'
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '1.985.89'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.05'
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5192'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.532'
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000122'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.775'
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05960'
$ws.Range('E51').Value = '  +0.85%  '
